$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Append the new daily rows (58 and 59) for 2025-11-29 (date serial 45990),
# one per station, matching the existing table's layout/columns.
$ws.Cells.Item(58, 1).Value = 45990
$ws.Cells.Item(58, 2).Value = "四方坪站"
$ws.Cells.Item(58, 3).Value = 8964.73
$ws.Cells.Item(58, 4).Value = 7994.32
$ws.Cells.Item(58, 5).Value = 2975.11
$ws.Cells.Item(58, 6).Value = 367

$ws.Cells.Item(59, 1).Value = 45990
$ws.Cells.Item(59, 2).Value = "高岭站"
$ws.Cells.Item(59, 3).Value = 5552.17
$ws.Cells.Item(59, 4).Value = 4905.66
$ws.Cells.Item(59, 5).Value = 1406.67
$ws.Cells.Item(59, 6).Value = 186

# Move/expand the visible selection to the new last cell, like Excel would
# leave it after data entry on the new last row.
$ws.Range("G59").Select()
